$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Split old row 27 (TestCase_E26) into two rows: keep row 27 for the
# "create a new watchlist" item, insert a new row 28 for TestCase_E27
# which carries the "name the watchlists" / "add description" items.

# First, insert a new row after row 27 so the remaining rows shift down.
$ws.Rows.Item(28).Insert()

# Row 27 (TestCase_E26 / OPQA-312 / create a new watchlist)
$ws.Cells.Item(27, 1).Value = "TestCase_E26"
$ws.Cells.Item(27, 2).Value = "OPQA-312"
$ws.Cells.Item(27, 3).Value = "Verify that user is able to create a new watchlist"
$ws.Cells.Item(27, 4).Value = "Y"
$ws.Cells.Item(27, 5).Value = "PASS"
$ws.Rows.Item(27).RowHeight = 15

# Row 28 (TestCase_E27 / OPQA-314 ||OPQA-317 / name watchlists + description)
$ws.Cells.Item(28, 1).Value = "TestCase_E27"
$ws.Cells.Item(28, 2).Value = "OPQA-314 ||OPQA-317"
$ws.Cells.Item(28, 3).Value = "Verify that user is able to name the watchlists||Verify that a user can add description to his watchlist"
$ws.Cells.Item(28, 4).Value = "Y"
$ws.Cells.Item(28, 5).Value = "PASS"
$ws.Rows.Item(28).RowHeight = 30

# Copy styles from row 27's original formatting (border + wrap) onto new cells
$ws.Range("A28").Style = $ws.Range("A27").Style
$ws.Range("B28").Style = $ws.Range("B27").Style
$ws.Range("C28").Style = $ws.Range("C27").Style
$ws.Range("D28").Style = $ws.Range("D27").Style
$ws.Range("E28").Style = $ws.Range("E27").Style

# Update selection/view state to match target
$ws.Range("E2:E27").Select()
$excel.ActiveWindow.ScrollRow = 37

$excel.ActiveWindow.WindowState = -4143
